$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.687.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.61%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.19%  "
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.70%  "
# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.04%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4794"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.05%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2907"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.11%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06578"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.23%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.69"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.21%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07786"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.05%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.55"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.40%  "
# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7418"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.47%  "
# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.869.02"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.52%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.196"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.91%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "280.99"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.47%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.688.88"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.62%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.52"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.05%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007623"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.33%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.152.29"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.79%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.297"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.30%  "
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.02%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.241"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.35%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.367"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.08%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.93"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.11%  "
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.32%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.971"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.70%  "
# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.01%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09978"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.38%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.520"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.14%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.365"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.14%  "
# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.26%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04787"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.28%  "
# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.19%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7060"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.81%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.12%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01877"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.76%  "
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.31%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.409"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.71%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.66"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.78%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4228"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.01%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.933"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.13%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8473"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.63%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.07%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.52"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.13%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.358"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.10%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.174"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.47%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "936.21"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.76%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.38"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.61%  "
# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.83%  "
